# Add 7T9J (actual) Omicron docking prediction results
# Fills in the HADDOCK metrics (columns K:Z) for the "Validation" rows
# (7T9J Chain C reference) on Sheet1 that previously only had the
# descriptive columns (A:J) populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column order (K..Z):
# HADDOCK score, HADDOCK score (+/-), Cluster size,
# RMSD from overall lowest-energy structure, RMSD (+/-),
# Van der Waals energy, Van der Waals energy (+/-),
# Electrostatic energy, Electrostatic energy (+/-),
# Desolvation energy, Desolvation energy (+/-),
# Restraints violation energy, Restraints violation energy (+/-),
# Buried Surface Area, Surface Area (+/-), Z-Score

$data = @{
    7  = @(-118.5, 23.8, 5, 1.2, 0.8, -76.4, 17.4, -169.6, 14.7, -20.8, 4.3, 125.8, 17.9, 2121.4, 257.3, -1.4)
    13 = @(-148.9, 1.1, 52, 0.6, 0.4, -94.1, 3.4, -191.1, 19.9, -30.6, 2.8, 140.4, 44.2, 2534.5, 33.4, -2.1)
    19 = @(-125.4, 12.6, 12, 0.8, 0.6, -81.9, 8.7, -223.5, 21.4, -8, 3.1, 92.5, 50.5, 2301.1, 54.4, -1.6)
    25 = @(-120, 4.2, 8, 16.5, 0.2, -67, 5.6, -271, 14.7, -17.2, 5.4, 183.8, 62.5, 1979.8, 99.6, -2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 11 + $i   # column K = 11
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Reflect the newly-populated columns in the frozen-pane scroll position
# and the active selection on the bottom-right pane.
$ws.Range("G17").Select() | Out-Null
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("G17").Select() | Out-Null
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("K26").Select() | Out-Null
